# The edit swaps the two species observation records currently stored in
# rows 23 and 24 (everything about the organism/position moves, while the
# shared locality/metadata columns that already match stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot current ("before") values for the columns that differ between
#     row 23 and row 24 ---
$a23 = $ws.Range("A23").Value()
$b23 = $ws.Range("B23").Value()
$e23 = $ws.Range("E23").Value()
$f23 = $ws.Range("F23").Value()
$g23 = $ws.Range("G23").Value()
$h23 = $ws.Range("H23").Value()
$q23 = $ws.Range("Q23").Value()
$r23 = $ws.Range("R23").Value()

$a24 = $ws.Range("A24").Value()
$b24 = $ws.Range("B24").Value()
$e24 = $ws.Range("E24").Value()
$f24 = $ws.Range("F24").Value()
$g24 = $ws.Range("G24").Value()
$h24 = $ws.Range("H24").Value()
$q24 = $ws.Range("Q24").Value()
$r24 = $ws.Range("R24").Value()
$ac24 = $ws.Range("AC24").Value()

# --- write row 23 with what used to be row 24's species/position data ---
$ws.Range("A23").Value = $a24
$ws.Range("B23").Value = $b24
$ws.Range("E23").Value = $e24
$ws.Range("F23").Value = $f24
$ws.Range("G23").Value = $g24
$ws.Range("H23").Value = $h24
$ws.Range("Q23").Value = $q24
$ws.Range("R23").Value = $r24
$ws.Range("AC23").Value = $ac24

# row 24 previously had empty "Ålder-Stadium/Kön/Aktivitet/Metod" cells
# (K24:N24) explicitly present; those move up to row 23 (touch them so the
# now-blank cells are materialised on row 23, matching default formatting).
$ws.Range("K23").Font.Bold = $false
$ws.Range("L23").Font.Bold = $false
$ws.Range("M23").Font.Bold = $false
$ws.Range("N23").Font.Bold = $false

# --- write row 24 with what used to be row 23's species/position data ---
$ws.Range("A24").Value = $a23
$ws.Range("B24").Value = $b23
$ws.Range("E24").Value = $e23
$ws.Range("F24").Value = $f23
$ws.Range("G24").Value = $g23
$ws.Range("H24").Value = $h23
$ws.Range("Q24").Value = $q23
$ws.Range("R24").Value = $r23

# the "Publik kommentar" and the (now unused) empty K:N cells leave row 24
$ws.Range("AC24").ClearContents()
$ws.Range("K24").ClearContents()
$ws.Range("L24").ClearContents()
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()
